# Updates cryptos list figures (prices/volume%) per Jun 4 2023 GitHub Actions run.
# Also swaps the ShibaInu/Avalanche row order (rows 16-17) to reflect new ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.275.00"
Set-TextValue "E2" "  +0.24%  "

Set-TextValue "D3" "1.908.87"
Set-TextValue "E3" "  +0.40%  "

Set-TextValue "E4" "  +0.32%  "

Set-TextValue "D5" "307.71"
Set-TextValue "E5" "  +0.01%  "

Set-TextValue "D7" "0.5373"
Set-TextValue "E7" "  +3.22%  "

Set-TextValue "D8" "0.3815"
Set-TextValue "E8" "  +1.19%  "

Set-TextValue "D9" "0.07299"
Set-TextValue "E9" "  +0.32%  "

Set-TextValue "D10" "22.07"
Set-TextValue "E10" "  +4.23%  "

Set-TextValue "D11" "0.9031"
Set-TextValue "E11" "  -0.02%  "

Set-TextValue "E12" "  -1.08%  "

Set-TextValue "E13" "  -0.63%  "

Set-TextValue "D14" "5.350"
Set-TextValue "E14" "  +1.40%  "

Set-TextValue "D15" "1.004"
Set-TextValue "E15" "  +0.39%  "

Set-TextValue "B16" "Avalanche"
Set-TextValue "C16" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D16" "14.83"
Set-TextValue "E16" "  +1.80%  "

Set-TextValue "B17" "ShibaInu"
Set-TextValue "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.000008647"
Set-TextValue "E17" "  +0.16%  "

Set-TextValue "E18" "  +0.27%  "

Set-TextValue "D19" "27.312.25"
Set-TextValue "E19" "  +0.30%  "

Set-TextValue "D20" "1.153.52"
Set-TextValue "E20" "  -39.35%  "

Set-TextValue "D21" "5.043"
Set-TextValue "E21" "  -0.93%  "

Set-TextValue "D22" "10.81"
Set-TextValue "E22" "  +1.46%  "

Set-TextValue "D23" "6.521"
Set-TextValue "E23" "  +1.50%  "

Set-TextValue "D24" "149.71"
Set-TextValue "E24" "  +1.80%  "

Set-TextValue "D25" "2.290"
Set-TextValue "E25" "  -1.34%  "

Set-TextValue "E26" "  +0.40%  "

Set-TextValue "D27" "1.746"
Set-TextValue "E27" "  -0.04%  "

Set-TextValue "D28" "116.85"
Set-TextValue "E28" "  +1.53%  "

Set-TextValue "D29" "4.835"
Set-TextValue "E29" "  -0.08%  "

Set-TextValue "D30" "4.803"
Set-TextValue "E30" "  -2.06%  "

Set-TextValue "D31" "0.09299"
Set-TextValue "E31" "  +0.55%  "

Set-TextValue "D32" "0.8353"
Set-TextValue "E32" "  +4.65%  "

Set-TextValue "E33" "  -0.04%  "

Set-TextValue "D34" "1.227"
Set-TextValue "E34" "  -1.01%  "

Set-TextValue "D35" "3.008"
Set-TextValue "E35" "  +2.25%  "

Set-TextValue "D36" "3.355"
Set-TextValue "E36" "  -2.20%  "

Set-TextValue "D37" "2.702"
Set-TextValue "E37" "  +4.01%  "

Set-TextValue "D38" "0.5771"
Set-TextValue "E38" "  +1.10%  "

Set-TextValue "D39" "0.02011"
Set-TextValue "E39" "  +0.37%  "

Set-TextValue "D40" "1.079"
Set-TextValue "E40" "  +0.16%  "

Set-TextValue "D41" "9.322"
Set-TextValue "E41" "  +3.18%  "

Set-TextValue "D42" "6.572"
Set-TextValue "E42" "  -0.18%  "

Set-TextValue "D43" "117.70"
Set-TextValue "E43" "  +0.86%  "

Set-TextValue "E44" "  +0.53%  "

Set-TextValue "D45" "0.4929"
Set-TextValue "E45" "  +1.49%  "

Set-TextValue "E46" "  +0.24%  "

Set-TextValue "D47" "10.16"
Set-TextValue "E47" "  +0.74%  "

Set-TextValue "D48" "1.640"
Set-TextValue "E48" "  +0.69%  "

Set-TextValue "D49" "38.51"
Set-TextValue "E49" "  +2.30%  "

Set-TextValue "D50" "0.06120"
Set-TextValue "E50" "  +2.80%  "

Set-TextValue "D51" "63.32"
Set-TextValue "E51" "  -0.92%  "

